# Apply the "updated GSC export data" change:
#  - Chart sheet: append 3 new daily rows (2025-11-09, 2025-11-10, 2025-11-11)
#  - Critical issues sheet: update "Not found (404)" page count, and
#    swap the "Discovered" / "Crawled - currently not indexed" rows with
#    their refreshed page counts.

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without Excel's COM layer
# re-interpreting date-like strings ("2025-11-09", etc.) as real dates.
# We stash the text behind a text formula in a scratch cell far away from
# any used range, copy/paste-special its computed (text) value into the
# destination, then clear the scratch cell again.
function Set-TextValue {
    param($ws, [int]$row, [int]$col, [string]$text)

    $scratch = $ws.Cells.Item(500, 500)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial(-4163) | Out-Null
    $scratch.ClearContents() | Out-Null
}

# ---------------------------------------------------------------------
# Chart sheet: new rows 37-39
# ---------------------------------------------------------------------
$chart = $wb.Worksheets.Item("Chart")

Set-TextValue $chart 37 1 "2025-11-09"
$chart.Cells.Item(37, 2).Value = 106
$chart.Cells.Item(37, 3).Value = 205
$chart.Cells.Item(37, 4).Value = 19

Set-TextValue $chart 38 1 "2025-11-10"
$chart.Cells.Item(38, 2).Value = 106
$chart.Cells.Item(38, 3).Value = 205
$chart.Cells.Item(38, 4).Value = 18

Set-TextValue $chart 39 1 "2025-11-11"
$chart.Cells.Item(39, 2).Value = 106
$chart.Cells.Item(39, 3).Value = 205
$chart.Cells.Item(39, 4).Value = 26

# ---------------------------------------------------------------------
# Critical issues sheet
# ---------------------------------------------------------------------
$critical = $wb.Worksheets.Item("Critical issues")

# Row 3: "Not found (404)" page count 8 -> 9
$critical.Cells.Item(3, 4).Value = 9

# Row 7 / Row 8: swap reasons, refreshed page counts
$critical.Cells.Item(7, 1).Value = "Crawled - currently not indexed"
$critical.Cells.Item(7, 4).Value = 7

$critical.Cells.Item(8, 1).Value = "Discovered - currently not indexed"
$critical.Cells.Item(8, 4).Value = 5
